$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 7.773964471275979
$ws.Cells.Item(2, 4).Value = 3.191376847968155
$ws.Cells.Item(2, 5).Value = 11.85124328904101
$ws.Cells.Item(2, 6).Value = 17.72866151351736
$ws.Cells.Item(2, 7).Value = 19.57760070707343
$ws.Cells.Item(2, 8).Value = 10.73842801354404
$ws.Cells.Item(2, 13).Value = 15.85493192387453
$ws.Cells.Item(2, 14).Value = 17.51857541418742
$ws.Cells.Item(2, 15).Value = 15.22314265709413
$ws.Cells.Item(3, 2).Value = 7.699799800734958
$ws.Cells.Item(3, 4).Value = 3.176492358786915
$ws.Cells.Item(3, 5).Value = 11.98720309686907
$ws.Cells.Item(3, 6).Value = 17.41093159309263
$ws.Cells.Item(3, 7).Value = 18.84099613348573
$ws.Cells.Item(3, 8).Value = 10.71732049147139
$ws.Cells.Item(3, 13).Value = 15.21398275592134
$ws.Cells.Item(3, 14).Value = 17.43140374411847
$ws.Cells.Item(3, 15).Value = 15.0646897611005
$ws.Cells.Item(4, 2).Value = 7.65569071955251
$ws.Cells.Item(4, 4).Value = 3.167747109919436
$ws.Cells.Item(4, 5).Value = 12.07665026904341
$ws.Cells.Item(4, 6).Value = 17.21898374929986
$ws.Cells.Item(4, 7).Value = 18.38104970430241
$ws.Cells.Item(4, 8).Value = 10.70688276961715
$ws.Cells.Item(4, 13).Value = 14.80740209072803
$ws.Cells.Item(4, 14).Value = 17.38013548785282
$ws.Cells.Item(4, 15).Value = 14.97138250067528
$ws.Cells.Item(5, 2).Value = 7.63809568132707
$ws.Cells.Item(5, 4).Value = 3.164285063018225
$ws.Cells.Item(5, 5).Value = 12.11458876781704
$ws.Cells.Item(5, 6).Value = 17.14167214614263
$ws.Cells.Item(5, 7).Value = 18.19204099147834
$ws.Cells.Item(5, 8).Value = 10.70326648453872
$ws.Cells.Item(5, 13).Value = 14.63866976506064
$ws.Cells.Item(5, 14).Value = 17.35982785094253
$ws.Cells.Item(5, 15).Value = 14.93440446794764
$ws.Cells.Item(6, 2).Value = 7.635197564504585
$ws.Cells.Item(6, 4).Value = 3.163716411632798
$ws.Cells.Item(6, 5).Value = 12.12097794577148
$ws.Cells.Item(6, 6).Value = 17.12889296341008
$ws.Cells.Item(6, 7).Value = 18.1605718907537
$ws.Cells.Item(6, 8).Value = 10.70270455978296
$ws.Cells.Item(6, 13).Value = 14.61047512102326
$ws.Cells.Item(6, 14).Value = 17.35649158801028
$ws.Cells.Item(6, 15).Value = 14.92832866719418
$ws.Cells.Item(7, 2).Value = 7.655451862316033
$ws.Cells.Item(7, 4).Value = 3.167700004265245
$ws.Cells.Item(7, 5).Value = 12.07715591203066
$ws.Cells.Item(7, 6).Value = 17.21793726181962
$ws.Cells.Item(7, 7).Value = 18.37850656358646
$ws.Cells.Item(7, 8).Value = 10.70683141588761
$ws.Cells.Item(7, 13).Value = 14.80513852134316
$ws.Cells.Item(7, 14).Value = 17.37985922241547
$ws.Cells.Item(7, 15).Value = 14.97087951306476
$ws.Cells.Item(8, 2).Value = 7.74810860835746
$ws.Cells.Item(8, 4).Value = 3.186164558427168
$ws.Cells.Item(8, 5).Value = 11.89687607134288
$ws.Cells.Item(8, 6).Value = 17.61852831001991
$ws.Cells.Item(8, 7).Value = 19.32541203265054
$ws.Cells.Item(8, 8).Value = 10.73062853309208
$ws.Cells.Item(8, 13).Value = 15.63676782551957
$ws.Cells.Item(8, 14).Value = 17.4880611371734
$ws.Cells.Item(8, 15).Value = 15.16770829606011
$ws.Cells.Item(9, 2).Value = 7.940205916979797
$ws.Cells.Item(9, 4).Value = 3.225376855139542
$ws.Cells.Item(9, 5).Value = 11.59127660455648
$ws.Cells.Item(9, 6).Value = 18.42375292565395
$ws.Cells.Item(9, 7).Value = 21.1076223322103
$ws.Cells.Item(9, 8).Value = 10.79715561059646
$ws.Cells.Item(9, 13).Value = 17.15540559663145
$ws.Cells.Item(9, 14).Value = 17.71734869796255
$ws.Cells.Item(9, 15).Value = 15.58326505363177
$ws.Cells.Item(10, 2).Value = 8.086393089064362
$ws.Cells.Item(10, 4).Value = 3.255847262261014
$ws.Cells.Item(10, 5).Value = 11.39680323723299
$ws.Cells.Item(10, 6).Value = 19.02040094467878
$ws.Cells.Item(10, 7).Value = 22.35521164358337
$ws.Cells.Item(10, 8).Value = 10.8578990124939
$ws.Cells.Item(10, 13).Value = 18.19278308360818
$ws.Cells.Item(10, 14).Value = 17.89513672870953
$ws.Cells.Item(10, 15).Value = 15.90376129555673
$ws.Cells.Item(11, 2).Value = 8.153716013097535
$ws.Cells.Item(11, 4).Value = 3.270031725518334
$ws.Cells.Item(11, 5).Value = 11.31504879103543
$ws.Cells.Item(11, 6).Value = 19.29150953180298
$ws.Cells.Item(11, 7).Value = 22.906521464844
$ws.Cells.Item(11, 8).Value = 10.88804575315348
$ws.Cells.Item(11, 13).Value = 18.6460360234555
$ws.Cells.Item(11, 14).Value = 17.97780011831567
$ws.Cells.Item(11, 15).Value = 16.05223190284211
$ws.Cells.Item(12, 2).Value = 8.179304197732563
$ws.Cells.Item(12, 4).Value = 3.275446179299991
$ws.Cells.Item(12, 5).Value = 11.2850726332545
$ws.Cells.Item(12, 6).Value = 19.394013982968
$ws.Cells.Item(12, 7).Value = 23.11275695449427
$ws.Cells.Item(12, 8).Value = 10.89981680292526
$ws.Cells.Item(12, 13).Value = 18.81487572973274
$ws.Cells.Item(12, 14).Value = 18.00933847288446
$ws.Cells.Item(12, 15).Value = 16.10878509219128
$ws.Cells.Item(13, 2).Value = 8.173789499065528
$ws.Cells.Item(13, 4).Value = 3.274278217247225
$ws.Cells.Item(13, 5).Value = 11.29148455752783
$ws.Cells.Item(13, 6).Value = 19.37194662551424
$ws.Cells.Item(13, 7).Value = 23.06845605457066
$ws.Cells.Item(13, 8).Value = 10.89726600994522
$ws.Cells.Item(13, 13).Value = 18.77863913972036
$ws.Cells.Item(13, 14).Value = 18.00253599250267
$ws.Cells.Item(13, 15).Value = 16.09659150907
$ws.Cells.Item(14, 2).Value = 8.155819410618916
$ws.Cells.Item(14, 4).Value = 3.27047633160772
$ws.Cells.Item(14, 5).Value = 11.31256284304768
$ws.Cells.Item(14, 6).Value = 19.29994644607143
$ws.Cells.Item(14, 7).Value = 22.92354041684207
$ws.Cells.Item(14, 8).Value = 10.8890070857585
$ws.Cells.Item(14, 13).Value = 18.65998322816193
$ws.Cells.Item(14, 14).Value = 17.98039020101016
$ws.Cells.Item(14, 15).Value = 16.05687825208295
$ws.Cells.Item(15, 2).Value = 8.144823812594774
$ws.Cells.Item(15, 4).Value = 3.268153078908956
$ws.Cells.Item(15, 5).Value = 11.32560237218237
$ws.Cells.Item(15, 6).Value = 19.25582031984175
$ws.Cells.Item(15, 7).Value = 22.83443994568714
$ws.Cells.Item(15, 8).Value = 10.88399430547807
$ws.Cells.Item(15, 13).Value = 18.58693564207881
$ws.Cells.Item(15, 14).Value = 17.96685529662133
$ws.Cells.Item(15, 15).Value = 16.03259418245858
$ws.Cells.Item(16, 2).Value = 8.082007882536793
$ws.Cells.Item(16, 4).Value = 3.254926492084018
$ws.Cells.Item(16, 5).Value = 11.40228261193562
$ws.Cells.Item(16, 6).Value = 19.00266727731055
$ws.Cells.Item(16, 7).Value = 22.31883803823119
$ws.Cells.Item(16, 8).Value = 10.85597890956121
$ws.Cells.Item(16, 13).Value = 18.16277634705056
$ws.Cells.Item(16, 14).Value = 17.88976861156639
$ws.Cells.Item(16, 15).Value = 15.89410760148898
$ws.Cells.Item(17, 2).Value = 8.04366588957298
$ws.Cells.Item(17, 4).Value = 3.246892829387218
$ws.Cells.Item(17, 5).Value = 11.45105472270605
$ws.Cells.Item(17, 6).Value = 18.84720417456817
$ws.Cells.Item(17, 7).Value = 21.99822648462984
$ws.Cells.Item(17, 8).Value = 10.83943194413986
$ws.Cells.Item(17, 13).Value = 17.8977031554553
$ws.Cells.Item(17, 14).Value = 17.84292090807557
$ws.Cells.Item(17, 15).Value = 15.80979730484905
$ws.Cells.Item(18, 2).Value = 8.021691340642361
$ws.Cells.Item(18, 4).Value = 3.24230271872824
$ws.Cells.Item(18, 5).Value = 11.47973788637209
$ws.Cells.Item(18, 6).Value = 18.75776329860993
$ws.Cells.Item(18, 7).Value = 21.81230230413272
$ws.Cells.Item(18, 8).Value = 10.83015164067591
$ws.Cells.Item(18, 13).Value = 17.74348984153923
$ws.Cells.Item(18, 14).Value = 17.81614482327965
$ws.Cells.Item(18, 15).Value = 15.76155931531226
$ws.Cells.Item(19, 2).Value = 8.014265385647711
$ws.Cells.Item(19, 4).Value = 3.240753950720825
$ws.Cells.Item(19, 5).Value = 11.48955732836687
$ws.Cells.Item(19, 6).Value = 18.72747982346402
$ws.Cells.Item(19, 7).Value = 21.74909757720819
$ws.Cells.Item(19, 8).Value = 10.82705039350984
$ws.Cells.Item(19, 13).Value = 17.69097903859151
$ws.Cells.Item(19, 14).Value = 17.80710865629823
$ws.Cells.Item(19, 15).Value = 15.74527221549625
$ws.Cells.Item(20, 2).Value = 8.047739481811854
$ws.Cells.Item(20, 4).Value = 3.247744879766326
$ws.Cells.Item(20, 5).Value = 11.44579744680711
$ws.Cells.Item(20, 6).Value = 18.86375666854734
$ws.Cells.Item(20, 7).Value = 22.03251465863914
$ws.Cells.Item(20, 8).Value = 10.84116890536857
$ws.Cells.Item(20, 13).Value = 17.92610267601242
$ws.Cells.Item(20, 14).Value = 17.84789054024965
$ws.Cells.Item(20, 15).Value = 15.81874627353357
$ws.Cells.Item(21, 2).Value = 8.161095285811111
$ws.Cells.Item(21, 4).Value = 3.271591895004359
$ws.Cells.Item(21, 5).Value = 11.30634483476283
$ws.Cells.Item(21, 6).Value = 19.32109982890854
$ws.Cells.Item(21, 7).Value = 22.96617588402351
$ws.Cells.Item(21, 8).Value = 10.89142334756346
$ws.Cells.Item(21, 13).Value = 18.69491209589894
$ws.Cells.Item(21, 14).Value = 17.98688873662197
$ws.Cells.Item(21, 15).Value = 16.0685344572899
$ws.Cells.Item(22, 2).Value = 8.23571955417863
$ws.Cells.Item(22, 4).Value = 3.287426938629131
$ws.Cells.Item(22, 5).Value = 11.22093960223017
$ws.Cells.Item(22, 6).Value = 19.6190270049933
$ws.Cells.Item(22, 7).Value = 23.56153442521493
$ws.Cells.Item(22, 8).Value = 10.92633416953673
$ws.Cells.Item(22, 13).Value = 19.1810266684278
$ws.Cells.Item(22, 14).Value = 18.07909480813804
$ws.Cells.Item(22, 15).Value = 16.23368740963037
$ws.Cells.Item(23, 2).Value = 8.195849577746568
$ws.Cells.Item(23, 4).Value = 3.278953769584496
$ws.Cells.Item(23, 5).Value = 11.26599134719457
$ws.Cells.Item(23, 6).Value = 19.46014303364234
$ws.Cells.Item(23, 7).Value = 23.24519823821506
$ws.Cells.Item(23, 8).Value = 10.90751476631666
$ws.Cells.Item(23, 13).Value = 18.92310755743518
$ws.Cells.Item(23, 14).Value = 18.02976513782708
$ws.Cells.Item(23, 15).Value = 16.14538599881832
$ws.Cells.Item(24, 2).Value = 8.045897594739166
$ws.Cells.Item(24, 4).Value = 3.247359578624197
$ws.Cells.Item(24, 5).Value = 11.44817225823869
$ws.Cells.Item(24, 6).Value = 18.85627347365812
$ws.Cells.Item(24, 7).Value = 22.0170179443028
$ws.Cells.Item(24, 8).Value = 10.84038289962439
$ws.Cells.Item(24, 13).Value = 17.91326891029734
$ws.Cells.Item(24, 14).Value = 17.84564327799179
$ws.Cells.Item(24, 15).Value = 15.81469971445676
$ws.Cells.Item(25, 2).Value = 7.887256811518137
$ws.Cells.Item(25, 4).Value = 3.214464507961031
$ws.Cells.Item(25, 5).Value = 11.66873468790714
$ws.Cells.Item(25, 6).Value = 18.20453811465739
$ws.Cells.Item(25, 7).Value = 20.63528238096971
$ws.Cells.Item(25, 8).Value = 10.777054644512
$ws.Cells.Item(25, 13).Value = 15.85493192387453
$ws.Cells.Item(25, 14).Value = 17.65360046530465
$ws.Cells.Item(25, 15).Value = 15.46796786507774
